$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.026.10'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '3.951.16'
$ws.Range("E3").Value = '  -1.78%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.674'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.76%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.744'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.166'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.91%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000316'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("D14").Value = '4.599.24'
$ws.Range("E14").Value = '  -1.17%  '
$ws.Range("D15").Value = '3.977.91'
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("E16").Value = '  +5.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.86'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("E19").Value = '  -0.77%  '
$ws.Range("D20").Value = '72.084.78'
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '430.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.69'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +12.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '95.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.21'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +17.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.44'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.90'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.04'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.59'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.91%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.131'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("B34").Value = 'Cosmos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '679.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '68.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.437'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.17%  '
$ws.Range("D38").Value = '0.0₃0847'
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.37'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.00%  '
$ws.Range("E40").Value = '  -3.76%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.27%  '
$ws.Range("E43").Value = '  -2.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0484'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.148'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.10%  '
$ws.Range("D51").Value = '2.782.25'
$ws.Range("E51").Value = '  +8.97%  '
